# Updated tests to run 30 trials
# Refresh the simulation-result values on the "Overall" and "Zones" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overall sheet - row 2 (Simulation 1 summary)
# ---------------------------------------------------------------
$overall = $wb.Worksheets.Item("Overall")

$overall.Range("B2").Value = 96
$overall.Range("C2").Value = 29
$overall.Range("D2").Value = 1.351266104806558
$overall.Range("E2").Value = 0.52222222222222203
$overall.Range("F2").Value = 1.5675384220024717
$overall.Range("G2").Value = 62
$overall.Range("H2").Value = 40
$overall.Range("I2").Value = 102
$overall.Range("J2").Value = 781
$overall.Range("K2").Value = 38

# ---------------------------------------------------------------
# Zones sheet - rows 2-14 (per zone summary)
# ---------------------------------------------------------------
$zones = $wb.Worksheets.Item("Zones")

# Row 2 - Zone 1
$zones.Range("B2").Value = 9
$zones.Range("C2").Value = 1
$zones.Range("D2").Value = 1.0089743589743587
$zones.Range("E2").Value = 0.26666666666666661
$zones.Range("F2").Value = 1.0708333333333331

# Row 3 - Zone 2
$zones.Range("B3").Value = 5
$zones.Range("C3").Value = 0
$zones.Range("D3").Value = 0.85277777777777775
$zones.Range("E3").Value = 0.39999999999999974
$zones.Range("F3").Value = 1.3055555555555556

# Row 4 - Zone 3
$zones.Range("B4").Value = 13
$zones.Range("C4").Value = 3
$zones.Range("D4").Value = 0.99750000000000016
$zones.Range("E4").Value = 0.43333333333333313
$zones.Range("F4").Value = 1.1385416666666668

# Row 5 - Zone 4
$zones.Range("B5").Value = 7
$zones.Range("C5").Value = 6
$zones.Range("D5").Value = 1.8766666666666665
$zones.Range("E5").Value = 0.15000000000000002
$zones.Range("F5").Value = 2

# Row 6 - Zone 5
$zones.Range("B6").Value = 5
$zones.Range("C6").Value = 3
$zones.Range("D6").Value = 1.4033333333333335
$zones.Range("E6").Value = 0.86666666666666625
$zones.Range("F6").Value = 1.462962962962963

# Row 7 - Zone 6 (E7 is a new cell in this edit)
$zones.Range("B7").Value = 9
$zones.Range("C7").Value = 3
$zones.Range("D7").Value = 1.9769230769230766
$zones.Range("E7").Value = 0.43749999999999956
$zones.Range("F7").Value = 2.661111111111111

# Row 8 - Zone 7
$zones.Range("B8").Value = 2
$zones.Range("C8").Value = 1
$zones.Range("D8").Value = 0.61111111111111072
$zones.Range("E8").Value = 0.4833333333333325
$zones.Range("F8").Value = 0.86666666666666714

# Row 9 - Zone 8 (E9 is removed in this edit)
$zones.Range("B9").Value = 5
$zones.Range("C9").Value = 0
$zones.Range("D9").Value = 2.1071428571428572
$zones.Range("E9").ClearContents()
$zones.Range("F9").Value = 2.1071428571428572

# Row 10 - Zone 9
$zones.Range("B10").Value = 4
$zones.Range("C10").Value = 3
$zones.Range("D10").Value = 1.3592592592592594
$zones.Range("E10").Value = 0.97500000000000009
$zones.Range("F10").Value = 1.4690476190476189

# Row 11 - Zone 10
$zones.Range("C11").Value = 0
$zones.Range("D11").Value = 1.1983333333333333
$zones.Range("E11").Value = 0.64583333333333359
$zones.Range("F11").Value = 1.5666666666666664

# Row 12 - Zone 11
$zones.Range("C12").Value = 5
$zones.Range("D12").Value = 1.5711111111111109
$zones.Range("E12").Value = 0.50833333333333275
$zones.Range("F12").Value = 1.7346153846153844

# Row 13 - Zone 12 (E13 is a new cell in this edit)
$zones.Range("B13").Value = 7
$zones.Range("C13").Value = 4
$zones.Range("D13").Value = 1.3589743589743588
$zones.Range("E13").Value = 0.61666666666666625
$zones.Range("F13").Value = 1.4208333333333332

# Row 14 - Zone 13
$zones.Range("B14").Value = 11
$zones.Range("D14").Value = 0.7333333333333335
$zones.Range("E14").Value = 0.51333333333333353
$zones.Range("F14").Value = 0.91666666666666685
